# lr! update sch/sth forms
# Rename the Liberia urine-filtration impact form from V2 to V2.1:
#  - "Select the council" (u_district label)  -> "Select the district"
#  - u_subdistrict field renamed to u_location, with label
#       "Select ward" -> "Select the location"
#  - repeat group id lr_u_202401_v2 -> lr_u_202401_v2_1
#  - settings sheet form_title / form_id bumped to the V2.1 / _v2_1 variants

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# Row 3: u_district label
$survey.Range("C3").Value = "Select the district"

# Row 4: u_subdistrict -> u_location
$survey.Range("B4").Value = "u_location"
$survey.Range("C4").Value = "Select the location"

# Row 9: begin repeat group name
$survey.Range("B9").Value = "lr_u_202401_v2_1"

# Settings sheet: form_title / form_id
$settings.Range("A2").Value = "(2024 Jan) - 4. SCH/STH – Urine Form V2.1"
$settings.Range("B2").Value = "lr_sch_sth_impact_202401_4_urine_v2_1"
